$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "1+4="
$t.Cell(1, 2).Range.Text = "3+38="
$t.Cell(1, 3).Range.Text = "17-9="
$t.Cell(1, 4).Range.Text = "33+4="
$t.Cell(1, 5).Range.Text = "30+47="
$t.Cell(2, 1).Range.Text = "95-11="
$t.Cell(2, 2).Range.Text = "18+9="
$t.Cell(2, 3).Range.Text = "30-12="
$t.Cell(2, 4).Range.Text = "27+25="
$t.Cell(2, 5).Range.Text = "5+92="
$t.Cell(3, 1).Range.Text = "65-62="
$t.Cell(3, 2).Range.Text = "75-3="
$t.Cell(3, 3).Range.Text = "6+93="
$t.Cell(3, 4).Range.Text = "67-14="
$t.Cell(3, 5).Range.Text = "35+38="
$t.Cell(4, 1).Range.Text = "46-34="
$t.Cell(4, 2).Range.Text = "80-3="
$t.Cell(4, 3).Range.Text = "46-35="
$t.Cell(4, 4).Range.Text = "66-18="
$t.Cell(4, 5).Range.Text = "5+61="
$t.Cell(5, 1).Range.Text = "60+3="
$t.Cell(5, 2).Range.Text = "22-15="
$t.Cell(5, 3).Range.Text = "20-8="
$t.Cell(5, 4).Range.Text = "58-38="
$t.Cell(5, 5).Range.Text = "18+25="
$t.Cell(6, 1).Range.Text = "6+83="
$t.Cell(6, 2).Range.Text = "7+58="
$t.Cell(6, 3).Range.Text = "11+42="
$t.Cell(6, 4).Range.Text = "86-69="
$t.Cell(6, 5).Range.Text = "94-34="
$t.Cell(7, 1).Range.Text = "83-41="
$t.Cell(7, 2).Range.Text = "82-34="
$t.Cell(7, 3).Range.Text = "58+21="
$t.Cell(7, 4).Range.Text = "55+27="
$t.Cell(7, 5).Range.Text = "9+50="
$t.Cell(8, 1).Range.Text = "22+44="
$t.Cell(8, 2).Range.Text = "90-34="
$t.Cell(8, 3).Range.Text = "79-44="
$t.Cell(8, 4).Range.Text = "86-11="
$t.Cell(8, 5).Range.Text = "11+29="
$t.Cell(9, 1).Range.Text = "52-49="
$t.Cell(9, 2).Range.Text = "6-4="
$t.Cell(9, 3).Range.Text = "17+75="
$t.Cell(9, 4).Range.Text = "53-38="
$t.Cell(9, 5).Range.Text = "92+6="
$t.Cell(10, 1).Range.Text = "20+34="
$t.Cell(10, 2).Range.Text = "84-2="
$t.Cell(10, 3).Range.Text = "85-49="
$t.Cell(10, 4).Range.Text = "89-84="
$t.Cell(10, 5).Range.Text = "28+69="
$t.Cell(11, 1).Range.Text = "21+67="
$t.Cell(11, 2).Range.Text = "59-50="
$t.Cell(11, 3).Range.Text = "17+55="
$t.Cell(11, 4).Range.Text = "99-86="
$t.Cell(11, 5).Range.Text = "18-5="
$t.Cell(12, 1).Range.Text = "90-2="
$t.Cell(12, 2).Range.Text = "98-15="
$t.Cell(12, 3).Range.Text = "20+24="
$t.Cell(12, 4).Range.Text = "0+71="
$t.Cell(12, 5).Range.Text = "2+72="
$t.Cell(13, 1).Range.Text = "69+25="
$t.Cell(13, 2).Range.Text = "21+54="
$t.Cell(13, 3).Range.Text = "67-31="
$t.Cell(13, 4).Range.Text = "31-5="
$t.Cell(13, 5).Range.Text = "93-14="
$t.Cell(14, 1).Range.Text = "71-42="
$t.Cell(14, 2).Range.Text = "99-87="
$t.Cell(14, 3).Range.Text = "39+22="
$t.Cell(14, 4).Range.Text = "15+79="
$t.Cell(14, 5).Range.Text = "74-14="
$t.Cell(15, 1).Range.Text = "31+38="
$t.Cell(15, 2).Range.Text = "33+41="
$t.Cell(15, 3).Range.Text = "91-18="
$t.Cell(15, 4).Range.Text = "88-36="
$t.Cell(15, 5).Range.Text = "42+21="
$t.Cell(16, 1).Range.Text = "29+60="
$t.Cell(16, 2).Range.Text = "55-37="
$t.Cell(16, 3).Range.Text = "54-0="
$t.Cell(16, 4).Range.Text = "44+20="
$t.Cell(16, 5).Range.Text = "80-57="
$t.Cell(17, 1).Range.Text = "6+1="
$t.Cell(17, 2).Range.Text = "3+65="
$t.Cell(17, 3).Range.Text = "84-61="
$t.Cell(17, 4).Range.Text = "44+19="
$t.Cell(17, 5).Range.Text = "44-31="
$t.Cell(18, 1).Range.Text = "92-63="
$t.Cell(18, 2).Range.Text = "72-35="
$t.Cell(18, 3).Range.Text = "83-80="
$t.Cell(18, 4).Range.Text = "64+15="
$t.Cell(18, 5).Range.Text = "78+11="
$t.Cell(19, 1).Range.Text = "20+13="
$t.Cell(19, 2).Range.Text = "89-32="
$t.Cell(19, 3).Range.Text = "79-66="
$t.Cell(19, 4).Range.Text = "79-11="
$t.Cell(19, 5).Range.Text = "43+14="
$t.Cell(20, 1).Range.Text = "6+76="
$t.Cell(20, 2).Range.Text = "6+73="
$t.Cell(20, 3).Range.Text = "16-2="
$t.Cell(20, 4).Range.Text = "87+2="
$t.Cell(20, 5).Range.Text = "50-49="
